$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.426.90"
$ws.Range("E2").Value = "  +0.53%  "

# Row 3
$ws.Range("D3").Value = "1.874.29"
$ws.Range("E3").Value = "  -0.24%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.92"
$ws.Range("E5").Value = "  +0.51%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  +0.00%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4825"
$ws.Range("E7").Value = "  -0.32%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2828"
$ws.Range("E8").Value = "  -1.75%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06550"
$ws.Range("E9").Value = "  -0.72%  "

# Row 10
$ws.Range("D10").Value = "1.884.26"
$ws.Range("E10").Value = "  +0.31%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07454"
$ws.Range("E11").Value = "  +2.29%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.46"
$ws.Range("E12").Value = "  -1.73%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.093"
$ws.Range("E13").Value = "  -2.03%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.19"
$ws.Range("E14").Value = "  +1.11%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6577"
$ws.Range("E15").Value = "  +0.34%  "

# Row 16
$ws.Range("D16").Value = "30.380.98"
$ws.Range("E16").Value = "  +0.56%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.33"
$ws.Range("E17").Value = "  -0.30%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9994"
$ws.Range("E18").Value = "  -0.02%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007649"
$ws.Range("E19").Value = "  -0.68%  "

# Row 20
$ws.Range("D20").Value = "2.116.37"
$ws.Range("E20").Value = "  +0.23%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.313"
$ws.Range("E21").Value = "  -0.14%  "

# Row 22
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "222.09"
$ws.Range("E22").Value = "  +13.66%  "

# Row 23
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9999"
$ws.Range("E23").Value = "  +0.03%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.195"
$ws.Range("E24").Value = "  +1.09%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.271"
$ws.Range("E25").Value = "  -0.38%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.52"
$ws.Range("E26").Value = "  +4.25%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.60"
$ws.Range("E27").Value = "  +2.90%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.975"
$ws.Range("E28").Value = "  +2.87%  "

# Row 29
$ws.Range("E29").Value = "  +1.14%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09405"
$ws.Range("E30").Value = "  +2.92%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.310"
$ws.Range("E31").Value = "  +0.83%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.019"
$ws.Range("E32").Value = "  -1.21%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05047"
$ws.Range("E33").Value = "  -1.31%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.217"
$ws.Range("E34").Value = "  +10.99%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7577"
$ws.Range("E35").Value = "  +5.34%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.702"
$ws.Range("E36").Value = "  -0.37%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01840"
$ws.Range("E37").Value = "  +2.23%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.625"
$ws.Range("E38").Value = "  -0.48%  "

# Row 39
$ws.Range("E39").Value = "  +2.08%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9082"
$ws.Range("E40").Value = "  -1.20%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.948"
$ws.Range("E41").Value = "  +2.17%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.79"
$ws.Range("E42").Value = "  +0.52%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4301"
$ws.Range("E43").Value = "  +0.35%  "

# Row 44
$ws.Range("E44").Value = "  +0.46%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.475"
$ws.Range("E45").Value = "  +0.93%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.12"
$ws.Range("E46").Value = "  -0.51%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1303"
$ws.Range("E47").Value = "  -1.29%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.481"
$ws.Range("E48").Value = "  +8.26%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.963"
$ws.Range("E49").Value = "  -2.11%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.27"
$ws.Range("E50").Value = "  +0.62%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3902"
$ws.Range("E51").Value = "  +1.98%  "
